# Bug fix: when there are no sales for a period, the report sheets should
# show a single placeholder row ("No hay ventas en este periodo") instead of
# product rows, and the header row loses its "Nombre_producto"/"descripcion"
# label (replaced by a plain 0) while the remaining header labels shift.
$wb = $excel.ActiveWorkbook

$noSalesText = "No hay ventas en este periodo"

# --- Sheets 1-3: "12 meses", "3 meses antes", "3 meses despues" ---
# Same layout: B1 header -> 0, row 3 (second product) removed entirely,
# row 2 collapses to just A2 (index) + B2 (message), C2:E2 cleared.
foreach ($idx in 1..3) {
    $ws = $wb.Worksheets.Item($idx)

    $ws.Range("B1").Value = 0

    $ws.Rows.Item(3).Delete()

    $ws.Range("B2").Value = $noSalesText
    $ws.Range("C2:E2").ClearContents()
}

# --- Sheet 4: "Comparativa" ---
# Header columns get relabeled/shifted (12mo_mensual column dropped,
# "descripcion" moves from B1 to E1), then same row collapse as above.
$ws4 = $wb.Worksheets.Item(4)

$ws4.Range("B1").Value = 0
$ws4.Range("C1").Value = "3mo_pre_mensual"
$ws4.Range("D1").Value = "3mo_post_mensual"
$ws4.Range("E1").Value = "descripcion"

$ws4.Rows.Item(3).Delete()

$ws4.Range("B2").Value = $noSalesText
$ws4.Range("C2:F2").ClearContents()
